$wb = $excel.ActiveWorkbook
$ws2 = $wb.Worksheets.Item("Tasks 02-04 to 02-11")
$ws3 = $wb.Worksheets.Item("Tasks 01-28 to 02-04")

# --- Sheet "Tasks 02-04 to 02-11" (physical sheet2.xml) edits ---

# Row 2: C2 and D2 get values, F2 status flips from TODO to Done
$ws2.Range("C2").Value = 0.5
$ws2.Range("D2").Value = 0
$ws2.Range("F11").Copy()
$ws2.Range("F2").PasteSpecial(-4122)
$ws2.Range("F2").Value = "Done"

# Row 3,4,5: Time spent reduced from 2 to 1
$ws2.Range("B3").Value = 1
$ws2.Range("B4").Value = 1
$ws2.Range("B5").Value = 1

# Row 7: Time spent increased from 3 to 5
$ws2.Range("B7").Value = 5

# Row 10: Status flips from In Progress back to TODO
$ws2.Range("F3").Copy()
$ws2.Range("F10").PasteSpecial(-4122)
$ws2.Range("F10").Value = "TODO"

$wb.Save()
